# Adding EB automation using appium
# Clear the "Pass" (J column) results for all testcase rows on the
# TestCases sheet, which removes the now-unused "Pass" shared string
# entry on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$rows = 2..24
foreach ($r in $rows) {
    $cell = $ws.Range("J$r")
    if ($cell.Value2 -eq "Pass") {
        $cell.ClearContents()
    }
}
